$d = $word.ActiveDocument
$LDQ = [char]8220
$RDQ = [char]8221

# ---------------------------------------------------------------------
# 1) "gamificado" paragraph: remove spell-check run splits, merge to one run
# ---------------------------------------------------------------------
$old1 = "En la pantalla del medio " + $LDQ + "Aprender-Jugar" + $RDQ + ", encontraremos las actividades basadas en aprendizaje gamificado y aprendizaje por medio de repetición espaciada. Dado el objetivo académico de este proyecto, esta pantalla será una de las ultimas de implementar. "
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $old1, 2) | Out-Null

# ---------------------------------------------------------------------
# 2) "Github Deskstop" / " o SourceTree" paragraph
#    Merge run1..run4 (up through "Deskstop") into one run, and
#    separately merge the " o " + "SourceTree" runs into another,
#    keeping them as two distinct runs (matches diff run boundaries).
# ---------------------------------------------------------------------
$old2a = "La aplicación se desarrollará en la última versión Android Studio, en su configuración con Java, el manejo de las versiones se hará en Git, en preferencia se usará Github Deskstop"
$d.Content.Find.Execute($old2a, $true, $false, $false, $false, $false, $true, 1, $false, $old2a, 2) | Out-Null

$old2b = " o SourceTree"
$d.Content.Find.Execute($old2b, $true, $false, $false, $false, $false, $true, 1, $false, $old2b, 2) | Out-Null

# ---------------------------------------------------------------------
# 3) "Lollipop" paragraph: merge run1..run3 (ending at "no se ") into
#    one run, without touching the following "hará" run.
# ---------------------------------------------------------------------
$old3 = "El software operará en dispositivos Android a partir de la versión 5.0 " + $LDQ + "Lollipop" + $RDQ + ", en cuanto a requerimientos de hardware, no se "
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $old3, 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Lecciones aprendidas / dificultades paragraph: full text replace
# ---------------------------------------------------------------------
$old4 = "Mencione las dificultades encontradas durante el desarrollo del proyecto. Además, haga alusión a las principales lecciones aprendidas durante el proceso."
$new4 = "Ya que se está desarrollando una aplicación móvil, se tuvo que aprender a usar las nuevas herramientas como lo son Android Studio, Gradle, entre otros. Además, dado que las pruebas de rendimiento requieren de grandes cantidades de datos, fue necesario aprender sobre el manejo e integración de bases de datos, lo que resultó un poco problemático en un comienzo."
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# ---------------------------------------------------------------------
# 5) Bibliography entry [1]: merge spell-check-split runs
# ---------------------------------------------------------------------
$old5 = "[1] Insor.gov.co. 2020. Preguntas Frecuentes " + [char]8211 + " INSOR " + [char]8211 + " INSORInstituto Nacional Para Sordos. "
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $old5, 2) | Out-Null

# ---------------------------------------------------------------------
# 6) Bibliography entry [2]: merge spell-check-split runs
# ---------------------------------------------------------------------
$old6 = "[2]2019. Plan Institucional 2019-2022. online] Available at: <http://www.insor.gov.co/home/descargar/plan_estrategico_NSOR_2019_2022V1.pdf [Accessed 8 September 2020]."
$d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $old6, 2) | Out-Null
